## Juno_TestDataSheet.xlsx - "Committing 4 object types"
##
## This commit appends a new automated-test-run row to the
## MasterTestDataSheet sheet (TC_CongaTemplates_MNCCollectionProcedureRecord)
## and updates the "last run" timestamp recorded on the CongaTemplateCreation
## sheet, then leaves the active selection on the next empty data row (B8).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. MasterTestDataSheet: populate row 6 with a new Conga Templates test
#    case, following the exact same layout as rows 2-5.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("MasterTestDataSheet")

$ws1.Range("A6").Value = "TC_CongaTemplates_MNCCollectionProcedureRecord"
$ws1.Range("B6").Value = "CongaTemplateCreation"
$ws1.Range("C6").Value = "https://test.salesforce.com/"
$ws1.Range("D6").Value = '..\\JunoAutomation\src\resources\\Juno_TestDataSheet.xlsx'
$ws1.Range("E6").Value = "GoogleChrome"

# Match the formatting used by the row directly above it (font/border/
# wrap/row-height) rather than the leftover blank-row formatting.
$ws1.Range("B5:E5").Copy()
$ws1.Range("B6:E6").PasteSpecial(-4122)
$ws1.Rows.Item(6).RowHeight = 26

# The last touched cell on this sheet is the next blank row.
$null = $ws1.Range("B8").Select()

# ---------------------------------------------------------------------
# 2. CongaTemplateCreation: record the latest automation run id.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CongaTemplateCreation")
$ws3.Range("E2").Value = "Auto_WedJan091104242019"

Write-Host "Committing 4 object types"
